{"js": "// Replace the division-problem text in each table cell per the commit's\n// regenerated problem set. Each old string is unique in the document, so a\n// scoped search-and-replace on context.document.body is unambiguous and\n// leaves every other run (formatting, date paragraph, etc.) untouched.\nconst replacements = [\n  [\"613\u00f77=\", \"281\u00f73=\"],\n  [\"720\u00f79=\", \"295\u00f77=\"],\n  [\"730\u00f72=\", \"568\u00f77=\"],\n  [\"383\u00f77=\", \"866\u00f79=\"],\n  [\"331\u00f73=\", \"377\u00f72=\"],\n  [\"789\u00f78=\", \"477\u00f72=\"],\n  [\"279\u00f74=\", \"922\u00f72=\"],\n  [\"383\u00f72=\", \"824\u00f73=\"],\n  [\"369\u00f73=\", \"709\u00f72=\"],\n  [\"760\u00f78=\", \"144\u00f78=\"],\n  [\"231\u00f72=\", \"642\u00f78=\"],\n  [\"397\u00f76=\", \"668\u00f76=\"],\n  [\"549\u00f77=\", \"369\u00f72=\"],\n  [\"857\u00f74=\", \"836\u00f77=\"],\n  [\"808\u00f76=\", \"930\u00f79=\"],\n  [\"512\u00f77=\", \"949\u00f72=\"],\n  [\"107\u00f78=\", \"889\u00f76=\"],\n  [\"490\u00f76=\", \"430\u00f73=\"],\n  [\"320\u00f78=\", \"616\u00f73=\"],\n  [\"372\u00f76=\", \"198\u00f76=\"],\n  [\"976\u00f75=\", \"752\u00f78=\"],\n  [\"708\u00f74=\", \"305\u00f74=\"],\n  [\"494\u00f73=\", \"719\u00f73=\"],\n  [\"722\u00f77=\", \"514\u00f74=\"],\n  [\"442\u00f79=\", \"713\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in each table cell per the commit's\n# regenerated problem set. Each old string is unique in the document, so a\n# document-wide Find/Replace (ReplaceAll) for each exact pair is unambiguous\n# and leaves every other run (formatting, date paragraph, etc.) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"613\u00f77=\", \"281\u00f73=\"),\n    @(\"720\u00f79=\", \"295\u00f77=\"),\n    @(\"730\u00f72=\", \"568\u00f77=\"),\n    @(\"383\u00f77=\", \"866\u00f79=\"),\n    @(\"331\u00f73=\", \"377\u00f72=\"),\n    @(\"789\u00f78=\", \"477\u00f72=\"),\n    @(\"279\u00f74=\", \"922\u00f72=\"),\n    @(\"383\u00f72=\", \"824\u00f73=\"),\n    @(\"369\u00f73=\", \"709\u00f72=\"),\n    @(\"760\u00f78=\", \"144\u00f78=\"),\n    @(\"231\u00f72=\", \"642\u00f78=\"),\n    @(\"397\u00f76=\", \"668\u00f76=\"),\n    @(\"549\u00f77=\", \"369\u00f72=\"),\n    @(\"857\u00f74=\", \"836\u00f77=\"),\n    @(\"808\u00f76=\", \"930\u00f79=\"),\n    @(\"512\u00f77=\", \"949\u00f72=\"),\n    @(\"107\u00f78=\", \"889\u00f76=\"),\n    @(\"490\u00f76=\", \"430\u00f73=\"),\n    @(\"320\u00f78=\", \"616\u00f73=\"),\n    @(\"372\u00f76=\", \"198\u00f76=\"),\n    @(\"976\u00f75=\", \"752\u00f78=\"),\n    @(\"708\u00f74=\", \"305\u00f74=\"),\n    @(\"494\u00f73=\", \"719\u00f73=\"),\n    @(\"722\u00f77=\", \"514\u00f74=\"),\n    @(\"442\u00f79=\", \"713\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
